$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column A (Datum): replace the free-text weekday/date strings with real
# Excel date values, formatted as dates (numFmtId 14 -> "mm-dd-yy" maps to
# the built-in short-date format).
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "10/27/2020"
$ws.Range("A6").Value = "10/29/2020"
$ws.Range("A7").Value = "11/05/2020"
$ws.Range("A8").Value = "11/10/2020"
$ws.Range("A9").Value = "11/12/2020"
$ws.Range("A10").Value = "11/17/2020"
$ws.Range("A11").Value = "11/19/2020"

$ws.Range("A5:A11").NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------------
# Column D (Tätigkeit): tweak the wording of the last logged entry.
# ---------------------------------------------------------------------------
$ws.Range("D11").Value = "Stresstest durchgeführt. Tieferes einlesen in rp.h library (Abtastraten , Buffer, …)"

# Hours worked / running total for the existing last row changed slightly.
$ws.Range("B11").Value = 3
$ws.Range("C11").Formula = "=B11+C10"

# ---------------------------------------------------------------------------
# Two new rows of logged hours.
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "11/24/2020"
$ws.Range("B12").Value = 1
$ws.Range("C12").Formula = "=B12+C11"
$ws.Range("D12").Value = "Zoom Meeting"

$ws.Range("A13").Value = "11/25/2020"
$ws.Range("B13").Value = 2
$ws.Range("C13").Formula = "=B13+C12"
$ws.Range("D13").Value = "recherchen buffer(auslesen und schreiben, ..), rp.h"

$ws.Range("A12:A13").NumberFormat = "mm-dd-yy"
$ws.Range("D12").WrapText = $true
$ws.Range("D13").WrapText = $true

# Row heights: rows 6,7,10,11,13 are "30" tall (two-line entries); row 11 was
# 45 before and is now 30 since the text got shorter.
$ws.Range("A11:E11").RowHeight = 30
$ws.Range("A13:E13").RowHeight = 30

# ---------------------------------------------------------------------------
# Column A width auto-fits to the much-shorter date strings now in use.
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 9.25

# ---------------------------------------------------------------------------
# Selection / view bookkeeping, mirroring where the author's cursor ended up.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("E12").Select()
